$wb = $excel.ActiveWorkbook

$monster = $wb.Worksheets.Item("monster")
$role = $wb.Worksheets.Item("role")

# --- monster sheet: add "level" (E) and "sub" (F) columns ---

# Header row (row 1) - same style as the existing headers (vertical-center, s="1")
# New shared strings must land in the order level(28), int(29), sub(30) to
# mirror the original commit, so touch E1/E2 (level/int) before F1 (sub).
$monster.Range("E1").Value = "level"
$monster.Range("E1").VerticalAlignment = -4108

# Type row (row 2) - plain "int" markers, no special style
$monster.Range("E2").Value = "int"

$monster.Range("F1").Value = "sub"
$monster.Range("F1").VerticalAlignment = -4108
$monster.Range("F2").Value = "int"

# Data rows 4-12: level / sub values
$monster.Range("E4").Value = 1
$monster.Range("F4").Value = 1

$monster.Range("E5").Value = 2
$monster.Range("F5").Value = 1

$monster.Range("E6").Value = 3
$monster.Range("F6").Value = 2

$monster.Range("E7").Value = 4
$monster.Range("F7").Value = 2

$monster.Range("E8").Value = 5
$monster.Range("F8").Value = 3

$monster.Range("E9").Value = 6
$monster.Range("F9").Value = 4

$monster.Range("E10").Value = 7
$monster.Range("F10").Value = 5

$monster.Range("E11").Value = 8
$monster.Range("F11").Value = 4

$monster.Range("E12").Value = 9
$monster.Range("F12").Value = 3

# New column D width (closest representable value to 23.25 given the
# engine's character-width quantization)
$monster.Columns.Item(4).ColumnWidth = 22.5

# --- switch the active tab from "role" to "monster" ---
[void]$role.Range("E17").Select()
[void]$monster.Activate()
[void]$monster.Range("D10").Select()
